# C1--C2-and-C3-PowerPoint.pptx edit
# 1) Re-apply the table on slide 16 with the built-in table style
#    {AC1E41F5-05BC-4396-8272-231214BB2731} instead of the custom
#    Table_0 style {5FDF55C3-A641-4369-8688-C302D92FA199}.
# 2) Swap the presentation's active theme colour scheme (currently the
#    "Integral" palette) for the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{AC1E41F5-05BC-4396-8272-231214BB2731}")
    }
}

# --- 2. Theme colours -------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (standard clrScheme order)
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rrggbb = $officeColors[$i - 1]
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $themeColors.Colors($i).RGB = $bgr
}
